$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add data for adult.csv row (row 4) - introduces the
# "(searched by simulated annealing algorithm)" shared string first.
$ws.Range("B4").Value = 0.78437754271765603
$ws.Range("C4").Value = 0.75866558177379895
$ws.Range("D4").Value = "(searched by simulated annealing algorithm)"

# Add new row 5 for heart.csv - introduces "heart.csv" shared string next.
$ws.Range("A5").Value = "heart.csv"
$ws.Range("B5").Value = 0.86885245901639296
$ws.Range("C5").Value = 0.52459016393442603
$ws.Range("D5").Value = "(searched by simulated annealing algorithm)"

# Update header row: insert new column E ("Clustered comonotonicity") and
# move the previous E1 ("Weighted avg of Naive Bayes & Comonotonicity") to F1.
$ws.Range("F1").Value = $ws.Range("E1").Value2
$ws.Range("E1").Value = "Clustered comonotonicity"

# Set column F width to match new column added in diff (stored width 38.5).
$ws.Columns.Item(6).ColumnWidth = 37.67

# Update selection to match the diff (active cell E3)
[void]$ws.Range("E3").Select()
